# Fruta / hortaliza, semanal
# Insert a new data row above the existing row 83 ("Granada" / "Wonderfull" /
# 04-06-2021 record) with an updated weekly reading, shifting all the
# following rows down by one (old row 141 becomes row 142).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 83; Excel shifts rows 83..141 down to 84..142
# and copies formatting (incl. the date number format on column D) from the
# row above into the freshly inserted row.
$ws.Rows.Item(83).Insert()

# The values that used to live in row 83 are now in row 84. Duplicate that
# record into the new row 83, then patch the handful of columns that differ
# for this new weekly entry.
$srcRow = 84
$dstRow = 83
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item($dstRow, $col).Value = $ws.Cells.Item($srcRow, $col).Value()
}

# Apply the updated values for this new record.
$ws.Range("D83").Value = 44762
$ws.Range("K83").Value = "Sin especificar"
$ws.Range("M83").Value = 250
$ws.Range("N83").Value = 13000
$ws.Range("O83").Value = 13000
$ws.Range("P83").Value = 13000
$ws.Range("Q83").Value = "$/bandeja 10 kilos"
$ws.Range("S83").Value = 1300
$ws.Range("T83").Value = 10
